# "Perfect best friend" exercise — add the potential-outcomes columns
# (D: treatment effect TE = C - B; E: treatment indicator D; F: observed
# outcome Y = E*C + (1-E)*B) for rows 2:13, then the summary block in
# row 16 (new "PI" header), row 17 (SDO / ATE / selection bias / PI-scaled
# HE bias), row 19/20 (ATE + SEL BIAS + HE BIAS label/check) and row 23
# ("ok").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Treatment indicator D (column E), rows 2:13 -----------------------
$treat = @{2=0; 3=0; 4=1; 5=1; 6=1; 7=0; 8=1; 9=0; 10=0; 11=0; 12=1; 13=0}
foreach ($r in 2..13) {
    $ws.Cells.Item($r, 5).Value = $treat[$r]
}

# --- TE = C - B (column D), rows 2:13 -----------------------------------
# D2 was entered on its own first, then D3:D13 filled as one block — so it
# keeps a standalone formula while D3:D13 becomes a shared-formula group
# (matches the authored ref="D3:D13" si="0" grouping).
$ws.Range("D2").Formula = "=C2-B2"
$ws.Range("D3:D13").Formula = "=C3-B3"

# --- Observed outcome Y = E*C + (1-E)*B (column F), rows 2:13 ----------
$ws.Range("F2").Formula = "=E2*C2 + (1 - E2)*B2"
$ws.Range("F3:F13").Formula = "=E3*C3 + (1 - E3)*B3"

# --- New summary header "PI" next to the existing row 16 labels --------
# Give G16 the same look as the rest of the row-16 label band (font/fill)
# by copying A16's formatting across before putting the text in.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("G16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G16").Value = "PI"

# --- Row 17: SDO, ATE, selection bias, HE-bias-adjusted term, PI -------
$ws.Range("A17").Formula = "=AVERAGE(F4:F6,F8,F12) - AVERAGE(F2:F3,F7,F9:F11,F13)"
$ws.Range("B17").Formula = "=AVERAGE(D2:D13)"
$ws.Range("C17").Formula = "=AVERAGE(B4:B6,B8,B12) - AVERAGE(B2:B3,B7,B9:B11,B13)"
$ws.Range("G17").Formula = "=SUM(E2:E13)/COUNT(E2:E13)"
$ws.Range("D17").Formula = "=(1-G17)*(AVERAGE(D4:D6,D8,D12)-AVERAGE(D2:D3,D7,D9:D11,D13))"

$ws.Range("A17:D17").NumberFormat = "0.00"
$ws.Range("G17").NumberFormat = "0.00"

# --- Row 19/20: ATE + SEL BIAS + HE BIAS check --------------------------
$ws.Range("B19").Value = "ATE+SEL BIAS + HE BIAS"
$ws.Range("B20").Formula = "=B17+C17+D17"
$ws.Range("B20").NumberFormat = "0.00"

# --- Row 23: "ok" ---------------------------------------------------------
$ws.Range("A23").Value = "ok"

$ws.Calculate()
